{"js": "// Add ISMRM presentations and learning german\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 1) Languages: add \"German (learning)\" after \"English (advance)\"\n// ---------------------------------------------------------------------\nconst englishPara = paragraphs.items.find(p => p.text === \"English (advance)\");\nif (!englishPara) {\n  throw new Error(\"Could not find 'English (advance)' paragraph\");\n}\n// insertParagraph copies the pPr (numbering + style) of the anchor paragraph,\n// which matches the Compact / ilvl=1 / numId=1004 list item used here.\nenglishPara.insertParagraph(\"German (learning)\", Word.InsertLocation.after);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Publications: rewrite the \"1st Annual IIBM...\" bullet to reference the\n//    2023 ISMRM & ISMRT Annual Meeting & Exhibition, and add a second\n//    bullet for the Bayesian Optimization presentation at the same\n//    conference.\n// ---------------------------------------------------------------------\nparagraphs.load(\"text\");\nawait context.sync();\nconst pubPara = paragraphs.items.find(p => p.text.indexOf(\"1st Annual IIBM PhD. Symposium\") !== -1);\nif (!pubPara) {\n  throw new Error(\"Could not find the publication paragraph\");\n}\n\nconst newParagraphsOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr>' +\n                '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1015\"/></w:numPr>' +\n              '</w:pPr>' +\n              '<w:r><w:t xml:space=\"preserve\">Intensity-based Deep Learning for SPION concentration estimation in MR</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">imaging,</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=\"preserve\">2023 ISMRM &amp; ISMRT Annual Meeting &amp; Exhibition</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">. Toronto Canada,</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">June 2023.</w:t></w:r>' +\n            '</w:p>' +\n            '<w:p>' +\n              '<w:pPr>' +\n                '<w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1015\"/></w:numPr>' +\n              '</w:pPr>' +\n              '<w:r><w:t xml:space=\"preserve\">Finding Optimal Regularization Parameter for Undersampled Reconstruction using</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">Bayesian Optimization,</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=\"preserve\">2023 ISMRM &amp; ISMRT Annual Meeting &amp; Exhibition</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">.</w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n              '<w:r><w:t xml:space=\"preserve\">Toronto Canada, June 2023.</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>';\n\npubPara.insertOoxml(newParagraphsOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Add ISMRM presentations and learning german\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Languages: add \"German (learning)\" after \"English (advance)\"\n# ---------------------------------------------------------------------\n$englishPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"English (advance)\") {\n        $englishPara = $p\n        break\n    }\n}\nif ($null -eq $englishPara) {\n    throw \"Could not find 'English (advance)' paragraph\"\n}\n\n# Insert a new paragraph mark after \"English (advance)\" and give it the\n# text \"German (learning)\". The new paragraph inherits the pPr (Compact\n# style + ilvl=1/numId=1004 numbering) of the paragraph it split from.\n$englishPara.Range.InsertParagraphAfter()\n$newLangPara = $englishPara.Next()\n$newLangPara.Range.Text = \"German (learning)\"\n\n# ---------------------------------------------------------------------\n# 2) Publications: rewrite the \"1st Annual IIBM...\" bullet to reference the\n#    2023 ISMRM & ISMRT Annual Meeting & Exhibition, and add a second\n#    bullet for the Bayesian Optimization presentation at the same\n#    conference.\n# ---------------------------------------------------------------------\n$pubPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.IndexOf(\"1st Annual IIBM PhD. Symposium\") -ge 0) {\n        $pubPara = $p\n        break\n    }\n}\nif ($null -eq $pubPara) {\n    throw \"Could not find the publication paragraph\"\n}\n\n$newParagraphsOoxml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1015\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Intensity-based Deep Learning for SPION concentration estimation in MR</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">imaging,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=\"preserve\">2023 ISMRM &amp; ISMRT Annual Meeting &amp; Exhibition</w:t></w:r><w:r><w:t xml:space=\"preserve\">. Toronto Canada,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">June 2023.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1015\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Finding Optimal Regularization Parameter for Undersampled Reconstruction using</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">Bayesian Optimization,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:iCs/><w:i/></w:rPr><w:t xml:space=\"preserve\">2023 ISMRM &amp; ISMRT Annual Meeting &amp; Exhibition</w:t></w:r><w:r><w:t xml:space=\"preserve\">.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">Toronto Canada, June 2023.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\n'@\n\n$pubPara.Range.InsertXML($newParagraphsOoxml)\n"}
